$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting of the existing last header cell (G1) into the new H1 header
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Set the new header label and the corresponding data value
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
